# "finestra incidenza 7gg centrata su ultimo g"
# Change the 7-day rolling window used for columns C ("somma mobile 7gg.")
# and D ("somma mobile 7gg. per 100mila abitanti") from a window CENTERED
# on the current row (r-3 .. r+3) to a window ENDING on the current row
# (r-6 .. r), i.e. centered on the last day of the window.
#
# Column B ("nuovi pos.") holds the daily new-cases counts and is left
# untouched; C is recomputed as the trailing 7-day sum of B, and D is that
# sum scaled to a rate per 100,000 inhabitants (population = 1199).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$firstDataRow = 2
$population = 1199
$factor = 100000 / $population

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $windowStart = $r - 6

    if ($windowStart -ge $firstDataRow) {
        $windowRange = $ws.Range("B" + $windowStart + ":B" + $r)
        $sum = $excel.WorksheetFunction.Sum($windowRange)

        # Only rewrite the cells when the new trailing-window sum actually
        # differs from what is already stored, so rows whose value is
        # unaffected by the re-centring keep their original literal/
        # precision untouched.
        $cCell = $ws.Cells.Item($r, 3)
        $existing = $cCell.Value()
        if (($existing -is [string]) -or ($existing -ne $sum)) {
            $cCell.Value = $sum
            $ws.Cells.Item($r, 4).Value = $sum * $factor
        }
    } else {
        # Fewer than 7 days of history are available yet: leave blank,
        # same as the first days of the sheet originally were. Only touch
        # cells that currently hold a real (non-blank) value so rows that
        # were already blank are left exactly as-is.
        $cCell = $ws.Cells.Item($r, 3)
        if (-not ($cCell.Value() -is [string])) {
            $cCell.ClearContents()
        }
        $dCell = $ws.Cells.Item($r, 4)
        if (-not ($dCell.Value() -is [string])) {
            $dCell.ClearContents()
        }
    }
}
